$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of columns B and C for rows 1 through 11
for ($r = 1; $r -le 11; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    $cVal = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 2).Value = $cVal
    $ws.Cells.Item($r, 3).Value = $bVal
}

# Update the active cell selection
$ws.Range("G10").Select()
